# Poland IV Liga - "Atualização de bases das ligas, do dia: 30-03-2024 às 19:32"
#
# The update re-orders a number of match rows (the rows keep their sequential
# index in column A, but all the match data in columns B:AC is shuffled around
# between a handful of rows that share the same kick-off date/time). We read
# the current values of the affected rows and write them back in the new
# order using native COM range Value2 gets/sets, which lets Excel itself take
# care of managing the shared string table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap rows 30 and 31 -----------------------------------------------
$row30 = $ws.Range("B30:AC30").Value2
$row31 = $ws.Range("B31:AC31").Value2
$ws.Range("B30:AC30").Value2 = $row31
$ws.Range("B31:AC31").Value2 = $row30

# --- Swap rows 41 and 42 -----------------------------------------------
$row41 = $ws.Range("B41:AC41").Value2
$row42 = $ws.Range("B42:AC42").Value2
$ws.Range("B41:AC41").Value2 = $row42
$ws.Range("B42:AC42").Value2 = $row41

# --- Swap rows 43 and 44 -----------------------------------------------
$row43 = $ws.Range("B43:AC43").Value2
$row44 = $ws.Range("B44:AC44").Value2
$ws.Range("B43:AC43").Value2 = $row44
$ws.Range("B44:AC44").Value2 = $row43

# --- Rotate rows 59, 60, 61 ---------------------------------------------
# new(59) = old(61); new(60) = old(59); new(61) = old(60)
$row59 = $ws.Range("B59:AC59").Value2
$row60 = $ws.Range("B60:AC60").Value2
$row61 = $ws.Range("B61:AC61").Value2
$ws.Range("B59:AC59").Value2 = $row61
$ws.Range("B60:AC60").Value2 = $row59
$ws.Range("B61:AC61").Value2 = $row60

# --- Swap rows 65 and 66 -----------------------------------------------
$row65 = $ws.Range("B65:AC65").Value2
$row66 = $ws.Range("B66:AC66").Value2
$ws.Range("B65:AC65").Value2 = $row66
$ws.Range("B66:AC66").Value2 = $row65

# --- Swap rows 73 and 74 -----------------------------------------------
$row73 = $ws.Range("B73:AC73").Value2
$row74 = $ws.Range("B74:AC74").Value2
$ws.Range("B73:AC73").Value2 = $row74
$ws.Range("B74:AC74").Value2 = $row73

# --- Swap rows 75 and 76 -----------------------------------------------
$row75 = $ws.Range("B75:AC75").Value2
$row76 = $ws.Range("B76:AC76").Value2
$ws.Range("B75:AC75").Value2 = $row76
$ws.Range("B76:AC76").Value2 = $row75

# --- Swap rows 86 and 87 -----------------------------------------------
$row86 = $ws.Range("B86:AC86").Value2
$row87 = $ws.Range("B87:AC87").Value2
$ws.Range("B86:AC86").Value2 = $row87
$ws.Range("B87:AC87").Value2 = $row86

# --- Rotate rows 108, 109, 110 ------------------------------------------
# new(108) = old(110); new(109) = old(108); new(110) = old(109)
$row108 = $ws.Range("B108:AC108").Value2
$row109 = $ws.Range("B109:AC109").Value2
$row110 = $ws.Range("B110:AC110").Value2
$ws.Range("B108:AC108").Value2 = $row110
$ws.Range("B109:AC109").Value2 = $row108
$ws.Range("B110:AC110").Value2 = $row109
